$d = $word.ActiveDocument
$vt = [char]11
$rsq = [char]8217   # right single quotation mark
$ldq = [char]8220   # left double quotation mark
$rdq = [char]8221   # right double quotation mark

# ---------------------------------------------------------------------------
# Locate the "Main Line" content paragraph (immediately follows the
# "Main Line" Heading1 paragraph) by index.
# ---------------------------------------------------------------------------
$mainLineHeadingIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs($i).Range.Text.TrimEnd([char]13, [char]7)
    if ($txt -eq "Main Line") {
        $mainLineHeadingIndex = $i
        break
    }
}
if ($mainLineHeadingIndex -eq 0) {
    throw "Could not locate 'Main Line' heading paragraph"
}
$contentIndex = $mainLineHeadingIndex + 1

# ---------------------------------------------------------------------------
# 1. Rewrite the "Main Line" content paragraph text.
# ---------------------------------------------------------------------------
$newMainLineText = (
    "The God of Nature guides the adorable creature through a dark cave system, to reach an exit. They find " +
    "an exit but enter the dead place (where there is a lot of fire). It becomes noticeable that a growth on the adorable creature" + $rsq + "s body is starting to glow slightly." +
    $vt + "The adorable creature traverses through the dead place, and eventually arrives in the Jungle (where the adorable creature originates from)." +
    $vt + "As they enter the jungle, the adorable creature seems to grow weaker. This grows further apparent when it get closer and closer to its home." +
    $vt + "The Adorable creature eventually arrives at its home where its mother is. As the mother caresses her child, as the camera moves further away from the scene. A cry of sadness can be heard (indicating there is something wrong (a.k.a. Adorable creature is dead). The Game ends. :" + $rsq + "("
)

$contentPara = $d.Paragraphs($contentIndex)
$contentRange = $d.Range($contentPara.Range.Start, $contentPara.Range.End)
$contentRange.Text = $newMainLineText

# ---------------------------------------------------------------------------
# 2. Insert a new "Notes:" Heading1 paragraph right after the content
#    paragraph.
# ---------------------------------------------------------------------------
$contentPara = $d.Paragraphs($contentIndex)
$contentPara.Range.InsertParagraphAfter()

$notesHeadingIndex = $contentIndex + 1
$notesHeadingPara = $d.Paragraphs($notesHeadingIndex)
$notesHeadingPara.Range.Text = "Notes:"
$notesHeadingPara.Style = "Heading1"

# ---------------------------------------------------------------------------
# 3. Insert the notes body paragraph right after the "Notes:" heading, and
#    make sure it uses the "Normal" style (Heading1's "next style"), since
#    InsertParagraphAfter otherwise carries the heading style forward.
# ---------------------------------------------------------------------------
$notesHeadingPara = $d.Paragraphs($notesHeadingIndex)
$notesHeadingPara.Range.InsertParagraphAfter()

$notesBodyIndex = $notesHeadingIndex + 1
$notesBodyPara = $d.Paragraphs($notesBodyIndex)
$notesBodyPara.Style = "Normal"

$newNotesText = (
    "Though the theming seems confusing due to the dead place seemingly being the second " + $ldq + "world" + $rdq + "/section of the game which could confuse players as it seems like the climax of the narrative." +
    $vt + "What we could do to continue with the theming whilst trying to reduce the ludonarrative dissonance is show that even though it seems like a climatic section, it has very little danger and is fairly peaceful. Maybe it could be possible to have no enemies in this level until right at the end as a way to prepare the players for the final section?" +
    $vt + "Also with the final section being a jungle like environment, maybe we could use animal nature as the main danger, to create that sense that maybe the Cave and the Dead place weren" + $rsq + "t that bad after all. " +
    $vt + "Finally we got to ensure the ending hurts. Because we are sick bastards and we need to control those lovely bags of emotions (i.e. the players) >:)"
)

$notesBodyPara = $d.Paragraphs($notesBodyIndex)
$notesBodyPara.Range.Text = $newNotesText

# ---------------------------------------------------------------------------
# 4. Relocate the hidden "_GoBack" bookmark to the very end of the notes
#    body paragraph (right before its paragraph mark).
# ---------------------------------------------------------------------------
$notesBodyPara = $d.Paragraphs($notesBodyIndex)
$endPos = $notesBodyPara.Range.End - 1
$bookmarkRange = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
